# Rename the worksheets
$wb = $excel.ActiveWorkbook

$wsOneVsRest = $wb.Worksheets.Item("OneVsRest")
$wsOneVsRest.Name = "BinaryRelevance"

$wsMultiClass = $wb.Worksheets.Item("Multi-Class Transformation")
$wsMultiClass.Name = "LabelPowerset"

# Update sheet view / selection state on the (renamed) BinaryRelevance sheet
$wsBinaryRelevance = $wb.Worksheets.Item("BinaryRelevance")
$wsBinaryRelevance.Activate()
$excel.ActiveWindow.ScrollRow = 10
$wsBinaryRelevance.Range("C33:L35").Select()

# Update selection state on the LabelPowerset sheet
$wsLabelPowerset = $wb.Worksheets.Item("LabelPowerset")
$wsLabelPowerset.Activate()
$wsLabelPowerset.Range("F43").Select()
